$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated standard-error values (in parentheses, under theta/lambda point
# estimates) for the new bootstrapping results.
# "proportion drinking" column (E) standard errors remain "(0.0)".

$ws.Range("C3").Value  = "(0.2)"
$ws.Range("D3").Value  = "(0.03)"

$ws.Range("C5").Value  = "(0.46)"
$ws.Range("D5").Value  = "(0.16)"

$ws.Range("C7").Value  = "(0.39)"
$ws.Range("D7").Value  = "(0.12)"

$ws.Range("C9").Value  = "(0.38)"
$ws.Range("D9").Value  = "(0.06)"

$ws.Range("C11").Value = "(0.7)"
$ws.Range("D11").Value = "(0.47)"

$ws.Range("C13").Value = "(0.64)"
$ws.Range("D13").Value = "(0.16)"

$ws.Range("C15").Value = "(0.0)"
$ws.Range("D15").Value = "(0.03)"
